$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-20 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-21 Sunday", 2) | Out-Null
$d.Content.Find.Execute("18+22=", $true, $false, $false, $false, $false, $true, 1, $false, "15+34=", 2) | Out-Null
$d.Content.Find.Execute("83-81=", $true, $false, $false, $false, $false, $true, 1, $false, "15-12=", 2) | Out-Null
$d.Content.Find.Execute("2+42=", $true, $false, $false, $false, $false, $true, 1, $false, "26+66=", 2) | Out-Null
$d.Content.Find.Execute("18+12=", $true, $false, $false, $false, $false, $true, 1, $false, "83-0=", 2) | Out-Null
$d.Content.Find.Execute("87+2=", $true, $false, $false, $false, $false, $true, 1, $false, "61+10=", 2) | Out-Null
$d.Content.Find.Execute("51-15=", $true, $false, $false, $false, $false, $true, 1, $false, "70-39=", 2) | Out-Null
$d.Content.Find.Execute("28-23=", $true, $false, $false, $false, $false, $true, 1, $false, "60+35=", 2) | Out-Null
$d.Content.Find.Execute("66-32=", $true, $false, $false, $false, $false, $true, 1, $false, "52-47=", 2) | Out-Null
$d.Content.Find.Execute("77-2=", $true, $false, $false, $false, $false, $true, 1, $false, "21+37=", 2) | Out-Null
$d.Content.Find.Execute("81-45=", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=", 2) | Out-Null
$d.Content.Find.Execute("6+53=", $true, $false, $false, $false, $false, $true, 1, $false, "35+3=", 2) | Out-Null
$d.Content.Find.Execute("69+10=", $true, $false, $false, $false, $false, $true, 1, $false, "67-21=", 2) | Out-Null
$d.Content.Find.Execute("97-42=", $true, $false, $false, $false, $false, $true, 1, $false, "88-54=", 2) | Out-Null
$d.Content.Find.Execute("59+14=", $true, $false, $false, $false, $false, $true, 1, $false, "25+2=", 2) | Out-Null
$d.Content.Find.Execute("40+55=", $true, $false, $false, $false, $false, $true, 1, $false, "11+17=", 2) | Out-Null
$d.Content.Find.Execute("80-28=", $true, $false, $false, $false, $false, $true, 1, $false, "94-40=", 2) | Out-Null
$d.Content.Find.Execute("78-35=", $true, $false, $false, $false, $false, $true, 1, $false, "7-5=", 2) | Out-Null
$d.Content.Find.Execute("44-28=", $true, $false, $false, $false, $false, $true, 1, $false, "41+31=", 2) | Out-Null
$d.Content.Find.Execute("19+26=", $true, $false, $false, $false, $false, $true, 1, $false, "57+9=", 2) | Out-Null
$d.Content.Find.Execute("22-4=", $true, $false, $false, $false, $false, $true, 1, $false, "69-18=", 2) | Out-Null
$d.Content.Find.Execute("12+24=", $true, $false, $false, $false, $false, $true, 1, $false, "85-84=", 2) | Out-Null
$d.Content.Find.Execute("52-6=", $true, $false, $false, $false, $false, $true, 1, $false, "94-29=", 2) | Out-Null
$d.Content.Find.Execute("69-57=", $true, $false, $false, $false, $false, $true, 1, $false, "88-38=", 2) | Out-Null
$d.Content.Find.Execute("70-14=", $true, $false, $false, $false, $false, $true, 1, $false, "41-11=", 2) | Out-Null
$d.Content.Find.Execute("11+65=", $true, $false, $false, $false, $false, $true, 1, $false, "67+27=", 2) | Out-Null
$d.Content.Find.Execute("83-9=", $true, $false, $false, $false, $false, $true, 1, $false, "88-60=", 2) | Out-Null
$d.Content.Find.Execute("67-37=", $true, $false, $false, $false, $false, $true, 1, $false, "76+10=", 2) | Out-Null
$d.Content.Find.Execute("86-3=", $true, $false, $false, $false, $false, $true, 1, $false, "32-26=", 2) | Out-Null
$d.Content.Find.Execute("64-31=", $true, $false, $false, $false, $false, $true, 1, $false, "49+45=", 2) | Out-Null
$d.Content.Find.Execute("55-38=", $true, $false, $false, $false, $false, $true, 1, $false, "75-31=", 2) | Out-Null
$d.Content.Find.Execute("27+21=", $true, $false, $false, $false, $false, $true, 1, $false, "67+0=", 2) | Out-Null
$d.Content.Find.Execute("92-10=", $true, $false, $false, $false, $false, $true, 1, $false, "67-32=", 2) | Out-Null
$d.Content.Find.Execute("64+3=", $true, $false, $false, $false, $false, $true, 1, $false, "62+26=", 2) | Out-Null
$d.Content.Find.Execute("97-5=", $true, $false, $false, $false, $false, $true, 1, $false, "90-41=", 2) | Out-Null
$d.Content.Find.Execute("49+4=", $true, $false, $false, $false, $false, $true, 1, $false, "22+74=", 2) | Out-Null
$d.Content.Find.Execute("27+15=", $true, $false, $false, $false, $false, $true, 1, $false, "77-51=", 2) | Out-Null
$d.Content.Find.Execute("5+59=", $true, $false, $false, $false, $false, $true, 1, $false, "76-39=", 2) | Out-Null
$d.Content.Find.Execute("28+34=", $true, $false, $false, $false, $false, $true, 1, $false, "6+20=", 2) | Out-Null
$d.Content.Find.Execute("25-9=", $true, $false, $false, $false, $false, $true, 1, $false, "53+20=", 2) | Out-Null
$d.Content.Find.Execute("77-28=", $true, $false, $false, $false, $false, $true, 1, $false, "45+41=", 2) | Out-Null
$d.Content.Find.Execute("51+16=", $true, $false, $false, $false, $false, $true, 1, $false, "4+14=", 2) | Out-Null
$d.Content.Find.Execute("85-41=", $true, $false, $false, $false, $false, $true, 1, $false, "83-55=", 2) | Out-Null
$d.Content.Find.Execute("98-96=", $true, $false, $false, $false, $false, $true, 1, $false, "8+58=", 2) | Out-Null
$d.Content.Find.Execute("72-41=", $true, $false, $false, $false, $false, $true, 1, $false, "27+45=", 2) | Out-Null
$d.Content.Find.Execute("90-72=", $true, $false, $false, $false, $false, $true, 1, $false, "65-27=", 2) | Out-Null
$d.Content.Find.Execute("27+35=", $true, $false, $false, $false, $false, $true, 1, $false, "8+56=", 2) | Out-Null
$d.Content.Find.Execute("56-46=", $true, $false, $false, $false, $false, $true, 1, $false, "98-57=", 2) | Out-Null
$d.Content.Find.Execute("54-41=", $true, $false, $false, $false, $false, $true, 1, $false, "69-42=", 2) | Out-Null
$d.Content.Find.Execute("98-76=", $true, $false, $false, $false, $false, $true, 1, $false, "42-15=", 2) | Out-Null
$d.Content.Find.Execute("54-2=", $true, $false, $false, $false, $false, $true, 1, $false, "22-3=", 2) | Out-Null
$d.Content.Find.Execute("61+22=", $true, $false, $false, $false, $false, $true, 1, $false, "4+37=", 2) | Out-Null
$d.Content.Find.Execute("79-38=", $true, $false, $false, $false, $false, $true, 1, $false, "43-0=", 2) | Out-Null
$d.Content.Find.Execute("65-46=", $true, $false, $false, $false, $false, $true, 1, $false, "34+39=", 2) | Out-Null
$d.Content.Find.Execute("9+0=", $true, $false, $false, $false, $false, $true, 1, $false, "33-12=", 2) | Out-Null
$d.Content.Find.Execute("79+18=", $true, $false, $false, $false, $false, $true, 1, $false, "96-11=", 2) | Out-Null
$d.Content.Find.Execute("36+46=", $true, $false, $false, $false, $false, $true, 1, $false, "24-20=", 2) | Out-Null
$d.Content.Find.Execute("68-64=", $true, $false, $false, $false, $false, $true, 1, $false, "0+26=", 2) | Out-Null
$d.Content.Find.Execute("25+41=", $true, $false, $false, $false, $false, $true, 1, $false, "16-8=", 2) | Out-Null
$d.Content.Find.Execute("58-36=", $true, $false, $false, $false, $false, $true, 1, $false, "97-32=", 2) | Out-Null
$d.Content.Find.Execute("20+3=", $true, $false, $false, $false, $false, $true, 1, $false, "9+81=", 2) | Out-Null
$d.Content.Find.Execute("94-14=", $true, $false, $false, $false, $false, $true, 1, $false, "6+62=", 2) | Out-Null
$d.Content.Find.Execute("71-43=", $true, $false, $false, $false, $false, $true, 1, $false, "20+15=", 2) | Out-Null
$d.Content.Find.Execute("61-2=", $true, $false, $false, $false, $false, $true, 1, $false, "39-20=", 2) | Out-Null
$d.Content.Find.Execute("15+37=", $true, $false, $false, $false, $false, $true, 1, $false, "33+40=", 2) | Out-Null
$d.Content.Find.Execute("26+25=", $true, $false, $false, $false, $false, $true, 1, $false, "45+39=", 2) | Out-Null
$d.Content.Find.Execute("37-16=", $true, $false, $false, $false, $false, $true, 1, $false, "52-26=", 2) | Out-Null
$d.Content.Find.Execute("53+34=", $true, $false, $false, $false, $false, $true, 1, $false, "23+34=", 2) | Out-Null
$d.Content.Find.Execute("26+39=", $true, $false, $false, $false, $false, $true, 1, $false, "71-16=", 2) | Out-Null
$d.Content.Find.Execute("17+24=", $true, $false, $false, $false, $false, $true, 1, $false, "78+17=", 2) | Out-Null
$d.Content.Find.Execute("71-36=", $true, $false, $false, $false, $false, $true, 1, $false, "94-91=", 2) | Out-Null
$d.Content.Find.Execute("33+19=", $true, $false, $false, $false, $false, $true, 1, $false, "14+10=", 2) | Out-Null
$d.Content.Find.Execute("5+28=", $true, $false, $false, $false, $false, $true, 1, $false, "23+32=", 2) | Out-Null
$d.Content.Find.Execute("39-28=", $true, $false, $false, $false, $false, $true, 1, $false, "79-44=", 2) | Out-Null
$d.Content.Find.Execute("47+20=", $true, $false, $false, $false, $false, $true, 1, $false, "24+37=", 2) | Out-Null
$d.Content.Find.Execute("32+43=", $true, $false, $false, $false, $false, $true, 1, $false, "18+55=", 2) | Out-Null
$d.Content.Find.Execute("97-9=", $true, $false, $false, $false, $false, $true, 1, $false, "8+0=", 2) | Out-Null
$d.Content.Find.Execute("37+15=", $true, $false, $false, $false, $false, $true, 1, $false, "57+7=", 2) | Out-Null
$d.Content.Find.Execute("87-36=", $true, $false, $false, $false, $false, $true, 1, $false, "3+77=", 2) | Out-Null
$d.Content.Find.Execute("51+19=", $true, $false, $false, $false, $false, $true, 1, $false, "17+81=", 2) | Out-Null
$d.Content.Find.Execute("14+72=", $true, $false, $false, $false, $false, $true, 1, $false, "61+20=", 2) | Out-Null
$d.Content.Find.Execute("22+63=", $true, $false, $false, $false, $false, $true, 1, $false, "12+11=", 2) | Out-Null
$d.Content.Find.Execute("51-16=", $true, $false, $false, $false, $false, $true, 1, $false, "3+72=", 2) | Out-Null
$d.Content.Find.Execute("81-79=", $true, $false, $false, $false, $false, $true, 1, $false, "44-38=", 2) | Out-Null
$d.Content.Find.Execute("74+7=", $true, $false, $false, $false, $false, $true, 1, $false, "2+68=", 2) | Out-Null
$d.Content.Find.Execute("55+36=", $true, $false, $false, $false, $false, $true, 1, $false, "7+63=", 2) | Out-Null
$d.Content.Find.Execute("84-0=", $true, $false, $false, $false, $false, $true, 1, $false, "58-39=", 2) | Out-Null
$d.Content.Find.Execute("36+4=", $true, $false, $false, $false, $false, $true, 1, $false, "40-8=", 2) | Out-Null
$d.Content.Find.Execute("54+25=", $true, $false, $false, $false, $false, $true, 1, $false, "97-75=", 2) | Out-Null
$d.Content.Find.Execute("6+31=", $true, $false, $false, $false, $false, $true, 1, $false, "98-88=", 2) | Out-Null
$d.Content.Find.Execute("52-50=", $true, $false, $false, $false, $false, $true, 1, $false, "24-13=", 2) | Out-Null
$d.Content.Find.Execute("63-1=", $true, $false, $false, $false, $false, $true, 1, $false, "81-74=", 2) | Out-Null
$d.Content.Find.Execute("56-7=", $true, $false, $false, $false, $false, $true, 1, $false, "26+32=", 2) | Out-Null
$d.Content.Find.Execute("56-35=", $true, $false, $false, $false, $false, $true, 1, $false, "73+21=", 2) | Out-Null
$d.Content.Find.Execute("29+62=", $true, $false, $false, $false, $false, $true, 1, $false, "78-24=", 2) | Out-Null
$d.Content.Find.Execute("53-13=", $true, $false, $false, $false, $false, $true, 1, $false, "58-13=", 2) | Out-Null
$d.Content.Find.Execute("71-34=", $true, $false, $false, $false, $false, $true, 1, $false, "2+6=", 2) | Out-Null
$d.Content.Find.Execute("16+25=", $true, $false, $false, $false, $false, $true, 1, $false, "47-31=", 2) | Out-Null
$d.Content.Find.Execute("15-0=", $true, $false, $false, $false, $false, $true, 1, $false, "28+20=", 2) | Out-Null
$d.Content.Find.Execute("23-20=", $true, $false, $false, $false, $false, $true, 1, $false, "49-41=", 2) | Out-Null
$d.Content.Find.Execute("70-4=", $true, $false, $false, $false, $false, $true, 1, $false, "38+27=", 2) | Out-Null
